# Fill in previously-empty data cells in the "graph_excel" template sheet.
# The sheet contains several small data tables (shareholding pattern,
# board composition, audit fees, director remuneration, borrowing limits,
# CSR, etc.) whose value cells were blank in the source template and are
# populated here with the sample figures from the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shareholding pattern (rows 4-8): year headers + percentages
$ws.Range("C4").Value = 2015
$ws.Range("D4").Value = 2014
$ws.Range("E4").Value = 2013
$ws.Range("F4").Value = 2012

$ws.Range("C5").Value = 43.66
$ws.Range("D5").Value = 42.17
$ws.Range("E5").Value = 42.23
$ws.Range("F5").Value = 42.1

$ws.Range("C6").Value = 20.42
$ws.Range("D6").Value = 22.53
$ws.Range("E6").Value = 20.81
$ws.Range("F6").Value = 14.08

$ws.Range("C7").Value = 14.48
$ws.Range("D7").Value = 11.26
$ws.Range("E7").Value = 11.89
$ws.Range("F7").Value = 16.16

$ws.Range("C8").Value = 21.44
$ws.Range("D8").Value = 24.04
$ws.Range("E8").Value = 25.07
$ws.Range("F8").Value = 27.66

# Directors liable to retire by rotation (rows 16-18)
$ws.Range("C16").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("C18").Value = 6

# ID / NID counts (rows 28-29)
$ws.Range("C28").Value = 0.63
$ws.Range("D28").Value = 0.38
$ws.Range("C29").Value = 0.75
$ws.Range("D29").Value = 0.25

# Executive compensation vs shareholder value (rows 38-42)
$ws.Range("B38").Value = 2011
$ws.Range("C38").Value = 1.0954
$ws.Range("D38").Value = 136.9

$ws.Range("B39").Value = 2012
$ws.Range("C39").Value = 1.3471
$ws.Range("D39").Value = 187.36

$ws.Range("B40").Value = 2013
$ws.Range("C40").Value = 1.6936
$ws.Range("D40").Value = 143.49

$ws.Range("B41").Value = 2014
$ws.Range("C41").Value = 1.6386
$ws.Range("D41").Value = 178.48

$ws.Range("B42").Value = 2015
$ws.Range("C42").Value = 1.546
$ws.Range("D42").Value = 225.53

# Variation in director's remuneration (rows 50-51)
$ws.Range("C50").Value = "NA"
$ws.Range("D50").Value = 0.77

$ws.Range("C51").Value = 1.55
$ws.Range("D51").Value = 0.09

# Dividend / EPS / Payout (rows 60-62)
$ws.Range("B60").Value = 2013
$ws.Range("C60").Value = 1.25
$ws.Range("D60").Value = 3.98
$ws.Range("E60").Value = 0.37

$ws.Range("B61").Value = 2014
$ws.Range("C61").Value = 1.25
$ws.Range("D61").Value = 3.88
$ws.Range("E61").Value = 0.37

$ws.Range("B62").Value = 2015
$ws.Range("C62").Value = 1.25
$ws.Range("D62").Value = 7.89
$ws.Range("E62").Value = 0.18

# Dividend / EPS / Payout summary (rows 72-74)
$ws.Range("B72").Value = 1.25
$ws.Range("C72").Value = 7.89
$ws.Range("D72").Value = 0.18

$ws.Range("B73").Value = 5
$ws.Range("C73").Value = 18.25
$ws.Range("D73").Value = 0.32

$ws.Range("B74").Value = 59.5
$ws.Range("C74").Value = 105.91
$ws.Range("D74").Value = 0.65

# Audit / Audit-Related / Non-Audit fees (rows 82-85)
$ws.Range("C82").Value = 2015
$ws.Range("D82").Value = 2014

$ws.Range("C83").Value = 50
$ws.Range("D83").Value = 0.225

$ws.Range("C84").Value = 0.06
$ws.Range("D84").Value = 54.25

$ws.Range("C85").Value = 0.285
$ws.Range("D85").Value = 0.338

# Audit fee table with new FY column headers (row 96) and values (97-99)
$ws.Range("C96").Value = "FY 12/13"
$ws.Range("D96").Value = "FY 13/14"
$ws.Range("E96").Value = "FY 14/15"

$ws.Range("C97").Value = 0.225
$ws.Range("D97").Value = 0.225
$ws.Range("E97").Value = 50

$ws.Range("C98").Value = 0.04
$ws.Range("D98").Value = 54.25
$ws.Range("E98").Value = 0.06

$ws.Range("C99").Value = 37.35
$ws.Range("D99").Value = 0.338
$ws.Range("E99").Value = 0.285

# Second "FY -1/" directors block (rows 532-537) is relabeled "FY -1/00"
# to distinguish it from the other director table (rows 110-115), which
# keeps its original "FY -1/" label.
$ws.Range("B532").Value = "FY -1/00"
$ws.Range("B533").Value = "FY -1/00"
$ws.Range("B534").Value = "FY -1/00"
$ws.Range("B535").Value = "FY -1/00"
$ws.Range("B536").Value = "FY -1/00"
$ws.Range("B537").Value = "FY -1/00"
